$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C28").Value = 8
$ws.Range("C30").Value = 8

$ws.Range("F21").Select() | Out-Null
